$d = $word.ActiveDocument

# 1) Merge "COMUNICATO UFFICIALE n. " + "1  04.05.205" into one run,
#    dropping the proofErr grammar markers around the second part.
$d.Content.Find.Execute("COMUNICATO UFFICIALE n. 1  04.05.205", $true, $false, $false, $false, $false, $true, 1, $false, "COMUNICATO UFFICIALE n. 1  04.05.205", 2) | Out-Null

# 2) Merge "S.GIOVANNI" + " EVAN." into one run, dropping proofErr markers.
$d.Content.Find.Execute("S.GIOVANNI EVAN.", $true, $false, $false, $false, $false, $true, 1, $false, "S.GIOVANNI EVAN.", 2) | Out-Null

# 3) Replace team name "REAL CIAMPINO" with "FORTITUDO POMEZIA".
$d.Content.Find.Execute("REAL CIAMPINO", $true, $false, $false, $false, $false, $true, 1, $false, "FORTITUDO POMEZIA", 2) | Out-Null
